$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Инвестиционные_проекты" ------------------------------------
# Update the data that changed, and append the new project row.
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A3").Value = "Ta ya togo se"
$ws1.Range("B3").Value = 3
$ws1.Range("C3").Value = 2

$ws1.Range("A4").Value = "Mafioznik"
$ws1.Range("B4").Value = 26
$ws1.Range("C4").Value = 22

# --- Sheet 2 (new): "Статистика" ------------------------------------------
# Add it after the last existing sheet so it lands at the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws2.Name = "Статистика"

$ws2.Range("A1").Value = "Отрасль"
$ws2.Range("B1").Value = "Кол-во проектов"
$ws2.Range("C1").Value = "Сумма инвестиций по отрасли"

# Match the bold / bordered / centered header look used on sheet 1.
$hdr = $ws2.Range("A1:C1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

$ws2.Range("A2").Value = "станкоинструментальная промышленность"
$ws2.Range("B2").Value = 3
$ws2.Range("C2").Value = 88

$ws2.Range("A3").Value = "сельскохозяйственное машиностроение"
$ws2.Range("B3").Value = 2
$ws2.Range("C3").Value = 6
